# VLAB #38436  fix error in writing to ss_summary.sso
#
# Adds a new "I/O_Change" column (G) to the change log, populates it for
# the existing rows, appends the new 3.30.08.03 change-log entry in row 38,
# and updates the "last updated" banner in F1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- New column G: header + per-row I/O_Change notes ----------------------
$ws.Range("G3").Value = "I/O_Change"
$ws.Range("G3").Font.Bold = $true

$g = @{}
$g[4]  = "Yes"
$g[5]  = "No"
$g[6]  = "No"
$g[7]  = "No"
$g[8]  = "Yes,  conditional in control file"
$g[9]  = "Yes,  conditional in control file"
$g[10] = "Yes,  conditional in data file"
$g[11] = "Yes, mandatory in data file"
$g[12] = "Yes, mandatory in control file"
$g[13] = "No"
$g[14] = "No"
$g[15] = "No"
$g[16] = "No"
$g[17] = "No"
$g[18] = "No"
$g[19] = "Yes, mandatory in control file"
$g[20] = "No"
$g[21] = "Yes, mandatory in control file"
$g[22] = "No"
$g[23] = "Yes,  conditional in control file"
$g[24] = "Yes, mandatory in forecast file"
$g[25] = "No"
$g[26] = "No"
$g[27] = "No"
$g[28] = "No"
$g[29] = "No"
$g[30] = "No"
$g[31] = "No"
$g[32] = "Yes,  conditional in control file"
$g[33] = "No"
$g[34] = "No"
$g[35] = "No"
$g[36] = "No"
$g[37] = "No"
$g[38] = "No"

for ($r = 4; $r -le 38; $r++) {
    $ws.Cells.Item($r, 7).Value = $g[$r]
}

# ---- Row 38: new 3.30.08.03 change-log entry ------------------------------
$ws.Cells.Item(38, 1).Value = 43007
$ws.Cells.Item(38, 2).Value = "3.30.08.03"
$ws.Cells.Item(38, 3).Value = "new"
$ws.Cells.Item(38, 4).Value = "putput"
$ws.Cells.Item(38, 6).Value = "fix sometimes fatal error in writing to the new summary output file:  ss_summary.sso"

# ---- F1 banner: bump the "last updated" note (added last so it lands as
#      the final new shared string, matching the source edit order) --------
$ws.Range("F1").Value = "2017-09-29 for 3.30.08.03"

# ---- View state: keep the active-cell selection in sync with the edit ----
$ws.Range("F31").Select()
